# Apply "Better workflow for Backbone" edits.

$wb = $excel.ActiveWorkbook

# --- demand sheet: replace timestamp labels with generic time indices ---
$demand = $wb.Worksheets.Item("demand")
$demand.Range("A3").Value = "t000001"
$demand.Range("A4").Value = "t000002"
$demand.Range("A5").Value = "t000003"
$demand.Range("A6").Value = "t000004"
$demand.Range("A7").Value = "t000005"

# --- node__unit__io sheet: add an "input" relation row for gas -> gas_turbine ---
$io = $wb.Worksheets.Item("node__unit__io")
$io.Range("A2").Value = "gas"
$io.Range("B2").Value = "gas_turbine"
$io.Range("C2").Value = "input"

# --- p_commodity_price sheet: keep only the first time step, reindex it ---
$price = $wb.Worksheets.Item("p_commodity_price")
$price.Range("A2").Value = "t000001"
$price.Rows.Item(3).Resize(4).Delete()
